$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("456501", "Bragg - Honey & Green Tea", "1", "30.99", "30.99"),
    @("456505", "Bragg - Ginger Lemon Honey", "1", "30.99", "30.99"),
    @("456503", "Bragg - Apple Cinn", "1", "30.99", "30.99"),
    @("462700", "Pobble Bubble Tea - Mango Red Dragon", "2", "50.60", "101.20")
)

$row = 4
foreach ($entry in $data) {
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = "'" + $entry[$col - 1]
    }
    $row++
}
